# Applies two edits to the "user manual" document:
#  1. "...A name and a correct NetID..." -> "...A name and a correct Yale NetID..."
#     (split into 3 runs: prefix | "Yale " | suffix)
#  2. "Close the app and open it, ..." -> "Close the library catalogue app
#     and reopen it, ..." (split into 5 runs: "Close the " | "library
#     catalogue " | "app and " | "re" | "open it, ...")
#
# Word normalises (re-coalesces) adjacent same-format runs whenever a new
# InsertBefore/InsertAfter happens later in the same paragraph, so all text
# insertions are done first, and the run splits (forced via a toggled
# character-formatting no-op) are applied afterwards, left to right.

$d = $word.ActiveDocument

function Split-RunAt($range) {
    # Force Word to keep `range` as its own run by nudging a character
    # formatting property and reverting it; this leaves the visible
    # formatting untouched but breaks the run at range's boundaries.
    $range.Font.Bold = $true
    $range.Font.Bold = $false
}

# ---------------------------------------------------------------------
# Edit 1: "A name and a correct NetID are required..." -> "...Yale NetID..."
# ---------------------------------------------------------------------

$r1 = $d.Content
$r1.Find.Execute("NetID are required to grant admin privileges.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r1.Collapse(1)
$r1.InsertBefore("Yale ")

$s1 = $d.Content
$s1.Find.Execute("Yale NetID are required to grant admin privileges.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$s1.SetRange($s1.Start, $s1.Start + 5)
Split-RunAt $s1

# ---------------------------------------------------------------------
# Edit 2: "Close the app and open it, ..." -> "Close the library catalogue
# app and reopen it, ..."
# ---------------------------------------------------------------------

$r2a = $d.Content
$r2a.Find.Execute("app and open it, and when the warning shows", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r2a.Collapse(1)
$r2a.InsertBefore("library catalogue ")

$r2b = $d.Content
$r2b.Find.Execute("open it, and when the warning shows", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r2b.Collapse(1)
$r2b.InsertBefore("re")

$s2a = $d.Content
$s2a.Find.Execute("library catalogue app and reopen", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$s2a.SetRange($s2a.Start, $s2a.Start + 19)
Split-RunAt $s2a

$s2b = $d.Content
$s2b.Find.Execute("app and reopen it, and when the warning shows", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$s2b.SetRange($s2b.Start + 8, $s2b.Start + 8 + 2)
Split-RunAt $s2b
